# "Merging of suites and updation of code"
#
# The ManageProducts "Input" sheet lists PSKU product codes in column B.
# Three of those product codes are being refreshed to newly generated
# codes (the workbook also grows its pool of placeholder product codes
# elsewhere, but the only observable grid change is these three cells).
#
#   B2 : prodisjY -> prodBEok
#   B5 : prodZedD -> prodWnkL
#   B8 : prodtJAD -> prodmfxh

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("B2").Value = "prodBEok"
$ws.Range("B5").Value = "prodWnkL"
$ws.Range("B8").Value = "prodmfxh"
